# Detailed Memory Scenes.docx - add theme tags to each body-area heading
# and relocate the trailing "_GoBack" bookmark to the blank paragraph that
# separates the "Torso Area" section from the "Arms Area" section.
#
# Commit message: "Updated story docs with the themes"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Append " | <Theme>" (bold, sz28) right after each heading's
#    existing text, as a *separate* run so the new text keeps its own
#    run-properties block (mirrors how Word splits runs when new text
#    is typed with formatting that happens to equal the neighbour's).
# ---------------------------------------------------------------------
function Add-ThemeTag {
    param($ParagraphIndex, $Suffix)

    $doc = $word.ActiveDocument
    $para = $doc.Paragraphs.Item($ParagraphIndex)
    $paraRange = $para.Range

    # Collapse to a point right before the paragraph mark (i.e. right
    # after the heading's existing text) and insert the new text there.
    $insertAt = $paraRange.End - 1
    $insertionPoint = $doc.Range($insertAt, $insertAt)
    $insertionPoint.InsertAfter($Suffix)

    # Re-select just the newly inserted characters and nudge a
    # formatting property (off/on) so the run commits as its own run
    # instead of silently merging into the preceding run.
    $newRange = $doc.Range($insertAt, $insertAt + $Suffix.Length)
    $newRange.Bold = 0
    $newRange.Bold = 1
}

Add-ThemeTag 3 " | Guilt"
Add-ThemeTag 11 " | Neglect"
Add-ThemeTag 19 " | Shame"

# ---------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the blank paragraph at the very
#    end of the document up to the blank paragraph right before the
#    "Arms Area - Aggression" heading.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$targetParagraph = $d.Paragraphs.Item(18)
$d.Bookmarks.Add("_GoBack", $targetParagraph.Range)
